# Auto-generated Word COM-interop script
# Applies the "Improve Web features document" edit to the
# Fonctionnalites WEB table: renumber/reword existing rows,
# fill in the "Signer un abonnement" row, append a new
# "Resilier un abonnement" row, and append a whole new
# "Gestion des reservations" section (3 rows).

$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

# Append the 4 new rows needed at the end of the table
for ($i = 0; $i -lt 4; $i++) {
    $t.Rows.Add() | Out-Null
}
$t = $d.Tables.Item(1)
Write-Host "Row count after add:" $t.Rows.Count

# Vertically merge column 1 of rows 2-6 ("1 - Gestion compte client" block)
$__top = $t.Cell(2, 1)
$__bottom = $t.Cell(6, 1)
$__top.Merge($__bottom) | Out-Null
$t = $d.Tables.Item(1)

# Row 2 / Col 1: "1 - Gestion compte client"
$__t = $d.Tables.Item(1)
$__cell = $__t.Cell(2, 1)
$__r = $__cell.Range
$__r.InsertXML('<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:jc w:val="center"/></w:pPr><w:r><w:t xml:space="preserve">1 </w:t></w:r><w:r><w:t>–</w:t></w:r><w:r><w:t xml:space="preserve"> Gestion compte client</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>') | Out-Null

# Row 2 / Col 2: "1) Creer un compte"
$__t = $d.Tables.Item(1)
$__cell = $__t.Cell(2, 2)
$__r = $__cell.Range
$__r.InsertXML('<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:jc w:val="center"/></w:pPr><w:r><w:t>1)</w:t></w:r><w:r><w:t xml:space="preserve"> Créer un compte</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>') | Out-Null

# Row 3 / Col 2: "2) Modifier ses informations"
$__t = $d.Tables.Item(1)
$__cell = $__t.Cell(3, 2)
$__r = $__cell.Range
$__r.InsertXML('<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:jc w:val="center"/></w:pPr><w:r><w:t>2</w:t></w:r><w:r><w:t>)</w:t></w:r><w:r><w:t xml:space="preserve"> Modifier ses informations</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>') | Out-Null

# Row 3 / Col 3: append "comme le mot de passe"
$__t = $d.Tables.Item(1)
$__cell = $__t.Cell(3, 3)
$__r = $__cell.Range
$__r.InsertXML('<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:jc w:val="center"/></w:pPr><w:r><w:t>L’utilisateur peut renseigner des informations supplémentaires comme le moyen de paiement ou modifier les informations déjà existantes comme le mot de passe</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>') | Out-Null

# Row 4 / Col 2: "3) Suppression de compte"
$__t = $d.Tables.Item(1)
$__cell = $__t.Cell(4, 2)
$__r = $__cell.Range
$__r.InsertXML('<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:jc w:val="center"/></w:pPr><w:r><w:t>3</w:t></w:r><w:r><w:t>)</w:t></w:r><w:r><w:t xml:space="preserve"> Suppression de compte</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>') | Out-Null

# Row 4 / Col 3: append "de la base de donnees"
$__t = $d.Tables.Item(1)
$__cell = $__t.Cell(4, 3)
$__r = $__cell.Range
$__r.InsertXML('<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:jc w:val="center"/></w:pPr><w:r><w:t>L’utilisateur ou l’administrateur doit pouvoir supprimer un compte de la base de données</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>') | Out-Null

# Row 5 / Col 2: "4) Signer un abonnement"
$__t = $d.Tables.Item(1)
$__cell = $__t.Cell(5, 2)
$__r = $__cell.Range
$__r.InsertXML('<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:jc w:val="center"/></w:pPr><w:r><w:t>4)</w:t></w:r><w:r><w:t xml:space="preserve"> Signer un abonnement</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>') | Out-Null

# Row 5 / Col 3: subscription explanation
$__t = $d.Tables.Item(1)
$__cell = $__t.Cell(5, 3)
$__r = $__cell.Range
$__r.InsertXML('<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:jc w:val="center"/></w:pPr><w:r><w:t>L’utilisateur peut choisir de souscrire un abonnement parmi ceux proposés, il doit cependant renseigner un moyen de paiement pour cela</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>') | Out-Null

# Row 6 / Col 2: "5) Resilier un abonnement"
$__t = $d.Tables.Item(1)
$__cell = $__t.Cell(6, 2)
$__r = $__cell.Range
$__r.InsertXML('<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:jc w:val="center"/></w:pPr><w:r><w:t>5</w:t></w:r><w:r><w:t>)</w:t></w:r><w:r><w:t xml:space="preserve"> Résilier un abonnement</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>') | Out-Null

# Row 6 / Col 3: cancellation explanation
$__t = $d.Tables.Item(1)
$__cell = $__t.Cell(6, 3)
$__r = $__cell.Range
$__r.InsertXML('<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:jc w:val="center"/></w:pPr><w:r><w:t>Si l’utilisateur à souscrit à un abonnement il peut l’annuler à tout moment</w:t></w:r><w:r><w:t xml:space="preserve"> ou lors de la suppression de son compte</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>') | Out-Null

# Row 7 / Col 1: "2 - Gestion des reservations"
$__t = $d.Tables.Item(1)
$__cell = $__t.Cell(7, 1)
$__r = $__cell.Range
$__r.InsertXML('<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:jc w:val="center"/></w:pPr><w:r><w:t xml:space="preserve">2 </w:t></w:r><w:r><w:t>–</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t>Gestion des réservations</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>') | Out-Null

# Row 7 / Col 2: "1) Effectuer une reservation"
$__t = $d.Tables.Item(1)
$__cell = $__t.Cell(7, 2)
$__r = $__cell.Range
$__r.InsertXML('<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:jc w:val="center"/></w:pPr><w:r><w:t>1) Effectuer une réservation</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>') | Out-Null

# Row 7 / Col 3: reservation explanation (+ page break)
$__t = $d.Tables.Item(1)
$__cell = $__t.Cell(7, 3)
$__r = $__cell.Range
$__r.InsertXML('<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:jc w:val="center"/></w:pPr><w:r><w:t xml:space="preserve">L’utilisateur peut réserver un </w:t></w:r><w:r><w:t xml:space="preserve">service en y indiquant le type de service puis d’autres informations selon le type de service (par exemple la fréquence des visites de proches âgés, ou le lieu et la date de réception d’un paquets), si l’utilisateur n’a pas d’abonnement </w:t></w:r><w:r><w:t xml:space="preserve">il devra avoir </w:t></w:r><w:r><w:lastRenderedPageBreak/><w:t xml:space="preserve">renseigner au préalable un moyen de paiement, </w:t></w:r><w:r><w:t xml:space="preserve">le compte de l’utilisateur sera débité </w:t></w:r><w:r><w:t>au moment de la réservation</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>') | Out-Null

# Row 8 / Col 2: "2) Modifier une reservation"
$__t = $d.Tables.Item(1)
$__cell = $__t.Cell(8, 2)
$__r = $__cell.Range
$__r.InsertXML('<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:jc w:val="center"/></w:pPr><w:r><w:t>2) Modifier une réservation</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>') | Out-Null

# Row 8 / Col 3: modify reservation explanation
$__t = $d.Tables.Item(1)
$__cell = $__t.Cell(8, 3)
$__r = $__cell.Range
$__r.InsertXML('<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:jc w:val="center"/></w:pPr><w:r><w:t>L’utilisateur peut modifier les informations d’une réservation</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>') | Out-Null

# Row 9 / Col 2: "3) Annuler une reservation"
$__t = $d.Tables.Item(1)
$__cell = $__t.Cell(9, 2)
$__r = $__cell.Range
$__r.InsertXML('<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:jc w:val="center"/></w:pPr><w:r><w:t>3) Annuler une réservation</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>') | Out-Null

# Row 9 / Col 3: cancel reservation explanation
$__t = $d.Tables.Item(1)
$__cell = $__t.Cell(9, 3)
$__r = $__cell.Range
$__r.InsertXML('<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:jc w:val="center"/></w:pPr><w:r><w:t>L’utilisateur peut annuler une réservation</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>') | Out-Null

# Re-create the _GoBack bookmark at its new location (end of row 7 / col 3)
$t = $d.Tables.Item(1)
$__cell73 = $t.Cell(7, 3)
$__r73 = $__cell73.Range
$__bmStart = $__r73.End - 1
$__bmRange = $d.Range($__bmStart, $__bmStart)
try {
    $d.Bookmarks.Add("_GoBack", $__bmRange) | Out-Null
} catch {
    Write-Host "Could not re-add _GoBack bookmark:" $_.Exception.Message
}

Write-Host "Done. Final row count:" $d.Tables.Item(1).Rows.Count

